# Applies the two changes described in the commit's XML diff:
#  1. Bump the cached "datetimeFigureOut" placeholder text from 7/16/13 to
#     7/17/13 on the Slide Master and on every Custom Layout (12 places).
#  2. On slide 1, split the run "puertos_disponibles" into "puertos" +
#     "EnUso" (two runs with identical formatting) inside the
#     "MensajesAgente" shape.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh (Slide Master + all Custom Layouts)
# ---------------------------------------------------------------------
$oldDate = "7/16/13"
$newDate = "7/17/13"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if (-not $isDatePlaceholder) { continue }

        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) "puertos_disponibles" -> "puertos" + "EnUso"
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($si = 1; $si -le $s.Shapes.Count; $si++) {
    $shp = $s.Shapes.Item($si)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }

    $tr = $shp.TextFrame.TextRange
    if ($tr.Text -notmatch "puertos_disponibles") { continue }

    $paraCount = $tr.Paragraphs(1, -1).Count
    for ($pi = 1; $pi -le 40; $pi++) {
        $para = $shp.TextFrame.TextRange.Paragraphs($pi, 1)
        if ($para.Text -eq "") { continue }
        if ($para.Text -notmatch "puertos_disponibles") { continue }

        for ($ri = 1; $ri -le 10; $ri++) {
            $run = $para.Runs($ri, 1)
            if ($run.Text -eq "puertos_disponibles") {
                $run.Text = "puertos"
                $run.InsertAfter("EnUso") | Out-Null
                break
            }
        }
    }
}

Write-Host "edit complete"
